$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "First_name"
$ws.Cells.Item(1, 2).Value = "Last_name"
$ws.Cells.Item(1, 3).Value = "Phone_number"

# --- Member rows 2-18 (re-sorted / re-ordered, with a few corrected spellings) ---
$firsts = @("Tarreessaa","Tolasaa","Ulfaataa","Baay``isaa","Caalaa","Qananiisaa","Walfaanaa","Lammii","Reggaasaa","Xurunaa","Fayyisaa","Silashii","Addunnaa","Katamaa","Biraanuu","Lalisee","Dabaree")
$lasts  = @("Wadaajoo","Calchisaa","Raggaasaa","Kabbadaa","Tolasaa","Biqilaa","Magarsaa","Diroo","Ballaxaa","Dabalee","Tolasaa","Dhabasaa","Nuurgii","Dhabasaa","Gaaddisaa","Magarsaa","Fayyeeraa")
$phones = @(913952050,923605989,912659004,922948389,910452943,921761067,912861288,922956646,913235855,948594041,911853155,985816078,913942964,925382373,910045632,954846351,912214364)

for ($i = 0; $i -lt $firsts.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $firsts[$i]
    $ws.Cells.Item($r, 2).Value = $lasts[$i]
    $ws.Cells.Item($r, 3).Value = $phones[$i]
    $ws.Cells.Item($r, 3).NumberFormat = "0000000000"
}

# --- Selection moved from F21 to F10 ---
$ws.Range("F10").Select()
